# Added New TestCase - TC_0002
#
# 1) Rename existing "Sheet1" to "Login".
# 2) Add a new worksheet "Create" right after "Login".
# 3) Populate "Create" with the same Uname/pwd columns as "Login" plus
#    companyName/firstName/lastName columns, for two new demo records
#    (TestLeaf / LeafTab test-company rows).
# 4) Carry over the header / data-row formatting (style) from "Login".
# 5) Leave "Login" selection at D4 (no longer the active tab) and make
#    "Create" the active tab with selection at G9.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Login"

# --- new worksheet, inserted right after "Login" -----------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Create"

# --- carry formatting over from the Login sheet -------------------------
# Header row (bold / filled / bordered) -> columns A:E
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)   # xlPasteFormats

# Data rows (bordered) -> columns A:E, rows 2:3
$ws1.Range("A2:B3").Copy()
$ws2.Range("A2:E3").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- header values --------------------------------------------------
$ws2.Range("A1").Value = "Uname"
$ws2.Range("B1").Value = "pwd"
$ws2.Range("C1").Value = "companyName"
$ws2.Range("D1").Value = "firstName"
$ws2.Range("E1").Value = "lastName"

# --- row 2 ------------------------------------------------------------
$ws2.Range("A2").Value = "DemoSalesManager"
$ws2.Range("B2").Value = "crmsfa"
$ws2.Range("C2").Value = "TestLeaf"
$ws2.Range("D2").Value = "TestFnameOne"
$ws2.Range("E2").Value = "TestLnameOne"

# --- row 3 ------------------------------------------------------------
$ws2.Range("A3").Value = "DemoCSR"
$ws2.Range("B3").Value = "crmsfa"
$ws2.Range("C3").Value = "LeafTab"
$ws2.Range("D3").Value = "TestFnameTwo"
$ws2.Range("E3").Value = "TestLnameTwo"

# --- column widths (best-effort autofit like Excel would do) -----------
$ws2.Columns.Item(1).EntireColumn.AutoFit()
$ws2.Columns.Item(3).EntireColumn.AutoFit()
$ws2.Columns.Item(4).EntireColumn.AutoFit()
$ws2.Columns.Item(5).EntireColumn.AutoFit()

# --- selections: Login left at D4, Create becomes the active tab at G9 -
[void]$ws1.Range("D4").Select()
[void]$ws2.Range("G9").Select()
